# Update "想去人数" (interest count) values in the "展览" and "全部类型" sheets
# to reflect newly generated output data.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 83
$wsExpo.Range("F4").Value = 255
$wsExpo.Range("F10").Value = 4409
$wsExpo.Range("F11").Value = 64

# Sheet "全部类型" (all types, aggregate of events)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 83
$wsAll.Range("F6").Value = 255
$wsAll.Range("F14").Value = 4409
$wsAll.Range("F15").Value = 64
